# add MLE by hand for all non-covariates
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update analyticDomain (O) and analyticRange (P) columns for rows 8-15
$ws.Range("O8").Value = "c(0,10)"
$ws.Range("P8").Value = "c(0,1)"

$ws.Range("O9").Value = "c(0,10)"
$ws.Range("P9").Value = "c(0,.5)"

$ws.Range("O10").Value = "c(0,20)"
$ws.Range("P10").Value = "c(0,.4)"

$ws.Range("O11").Value = "c(0,20)"
$ws.Range("P11").Value = "c(0,.5)"

$ws.Range("O12").Value = "c(0,30)"
$ws.Range("P12").Value = "c(0,.5)"

$ws.Range("O13").Value = "c(0,5)"
$ws.Range("P13").Value = "c(0,1.5)"

$ws.Range("O14").Value = "c(0,5)"
$ws.Range("P14").Value = "c(0,1.5)"

$ws.Range("O15").Value = "c(0,5)"
$ws.Range("P15").Value = "c(0,1.5)"

# Adjust column widths: split old N:P (14-16) block of width 16 into
# N=16, O=10, P=9.43
$ws.Columns.Item(15).ColumnWidth = 10
$ws.Columns.Item(16).ColumnWidth = 9.42578125

# Update the active selection to O13:O15
$ws.Range("O13:O15").Select()
